$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the refreshed coin prices and 1h volume-change percentages produced by
# the scheduled GitHub Actions job. Values that look like plain numbers must be
# forced to stay text (as in the source data) without leaving a lingering custom
# number format on the cell, so we flip to text format, assign, then restore the
# cell style back to Normal.
$ws.Range('D2').Value = '58.326.83'
$ws.Range('E2').Value = '  -3.03%  '
$ws.Range('D3').Value = '2.297.77'
$ws.Range('E3').Value = '  -4.83%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '545.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.79'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.95%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.570'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.84%  '
$ws.Range('D9').Value = '2.293.88'
$ws.Range('E9').Value = '  -4.91%  '
$ws.Range('E10').Value = '  -4.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.49'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.19%  '
$ws.Range('E12').Value = '  +0.91%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.332'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.64'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.72%  '
$ws.Range('D15').Value = '2.705.06'
$ws.Range('E15').Value = '  -5.00%  '
$ws.Range('D16').Value = '58.258.28'
$ws.Range('E16').Value = '  -3.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000132'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.36%  '
$ws.Range('D18').Value = '2.257.74'
$ws.Range('E18').Value = '  -6.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.60'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.28'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '312.93'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.44'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.90%  '
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.02'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.87%  '
$ws.Range('E25').Value = '  -4.16%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.82%  '
$ws.Range('E28').Value = '  -6.44%  '
$ws.Range('E29').Value = '  -3.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.23'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.14%  '
$ws.Range('D31').Value = '0.0₃0719'
$ws.Range('E31').Value = '  -6.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.09'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.74'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.380'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.73'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.80%  '
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.24'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.93'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '38.05'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.26%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.50'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '291.18'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -10.18%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '140.28'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.42'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0951'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0500'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.32%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.553'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.31'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -7.59%  '
$ws.Range('E49').Value = '  -3.60%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.53'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.75%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '10.99'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.57%  '
